# Edit Review_332.docx per commit diff:
# - Update date in paragraph 1 (29.10.24 -> 28.10.24) and replace subtitle after the line break
# - Replace body text of paragraphs 2-6
# - Remove the "חייב להגיד..." paragraph entirely
# - Replace the arxiv URL paragraph with the new link

$d = $word.ActiveDocument

# --- Paragraph 1: two runs of text joined by a <w:br/> ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range.Find
$r1.Execute("-29.10.24:", $true, $false, $false, $false, $false, $true, 1, $false, "-28.10.24:", 2) | Out-Null

$p1b = $d.Paragraphs.Item(1)
$r1b = $p1b.Range.Find
$r1b.Execute('Global Lyapunov functions: a long-standing open problem in mathematics, with symbolic transformers', $true, $false, $false, $false, $false, $true, 1, $false, 'HEAVY-TAILED DIFFUSION MODELS', 2) | Out-Null

# --- Paragraph 2 ---
$d.Paragraphs.Item(2).Range.Text = 'המאמר עם השם הקצר הזה משך את עיניי כי יש לי חיבה גם למודלי דיפוזיה גנרטיביים וגם להתפלגויות בעלות תכונות מעניינות למשל זנבות כבדים. בגדול התפלגות נקראת בעלת זנב כבד או ארוך כאשר התפלגות לזנב שלה (כלומר המסה ההסתברותית מימין לנקודה) מקשל (הסתברות) הינה גבוה יותר מאשר להתפלגות מעריכית. נשמע קצת מסובך אבל במילים פשוטות ניתן להגיד כי להתפלגויות בעלות זנב כבד(HT) יש יותר מסה בקצוות.'

# --- Paragraph 3 ---
$d.Paragraphs.Item(3).Range.Text = 'למשל התפלגות נורמלית אינה בעלת זנבות כבדים והתפלגות סטודנט t וגם התפלגות קושי הן כן. אוקיי, למה אני בכלל מדבר על זה? הסיבה היא די פשוטה - ההנחה שנוכל להניח התפלגות גאוסית על כל סוג של דאטה אינה נכונה. יש סוגי דאטה שלא ניתן לאפיין אותם בצורה טוב עם התפלגות בעלות זנבות קלים. עקב גם אנו נתקשה לגנרט דאטה מהתפלגויות אלו אם נמדל אותו (הדאטה) עם מודלי הבנויים על הנחות גאוסיות גם אם המודלים האלו הם בעלי expressiveness גבוהה כמו מודלי הדיפוזיה. עדיין יהיה מאוד בעייתי ליצור באמצעותם דאטה בעלת התפלגות HT במיוחד בקצוות ההתפלגות.'

# --- Paragraph 4 ---
$d.Paragraphs.Item(4).Range.Text = 'אז המאמר, שהוא אחד הכבדים ביותר מתמטית מאלו שראיתי לאחרונה, מציע להחליף את התפלגויות גאוסיות שיש לנו במודלי דיפוזיה בהתפלגות סטודנט שהיא התפלגות HT. כלומר כל מה שהיה בעלת התפלגות גאוסית במודל דיפוזיה מקורי יהיה מהתפלגות t. דרך אגב אחד הפרמטרים של התפלגות t (שהיא כמובן וקטורית עבור מודלים אלו כי אנו רוצים לגנרט דאטה בעלת מימדים רבים) שהוא שולט ב״כבדות הזנב״ שלה וכאשר היא שואפת לאינסוף אנו מקבלים את ההתפלגות הגאוסית האהובה עלינו. כלומר המודלים המוצעים במאמר הם הכללה של מודלי דיפוזיה גאוסיים שאנו מכירים ואוהבים.'

# --- Paragraph 5 ---
$d.Paragraphs.Item(5).Range.Text = 'כמובן שלא מספיק סתם להחליף התפלגות גאוסית במודל דיפוזיה בהתפלגות t - זה דורש להגדיר לא מעט התפלגויות מותנות הנדרשות לנו להגדרת הלמידה של תהליך denoising. זה די לא טריוויאלי אבל העקרון נשאר דומה -מאמנים את המודל להסיר רעש (שהוא מפולג עם t) באופן הדרגתי. במקום KL divergence המוכר לנו ממודלי דיפוזיה המחברים משתמשים ב-γ-Power divergence כדי למדוד מרחק בין ההתפלגות הדאטה אחרי הסרת רעש לזה של הדאטה האמיתי (לכל איטרציה).'

# --- Paragraph 6 ---
$d.Paragraphs.Item(6).Range.Text = 'גם תהליך הגנרוט מוגדר דומה עקרונות למודלי דיפוזיה גאוסיים אבל כמובן כל ה-hyperparameters מותאמים להתפלגות t. יש גם רפרמטריזציות שאנו כה אוהבים במודלי דיפוזיה, ייצוג באמצעות משוואות דיפרנציאליות חלקיות, גם באמצעות טכניקה חדשה הנקראת flow matching (הבונה מסלול מיטבי בין ההתפלגות ההתחלתית והתפלגות הדאטה). כאמור מאמר די כבד מתמטית ומקווה שהצלחתי להסביר לכם את העקרונות לפחות. '

# --- Paragraph 7 ("חייב להגיד...") is deleted entirely ---
$d.Paragraphs.Item(7).Range.Delete()

# --- Old paragraph 8 (now paragraph 7, the URL) gets new link text ---
$d.Paragraphs.Item(7).Range.Text = 'https://arxiv.org/pdf/2410.14171'

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
